# Add a new "SAVE" synonym row into the BLUESKY-COMMAND-TABLE sheet.
# It is inserted right above the existing "START" row (old row 108),
# pushing START..VRESOMETH down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 108 (shifts 108..113 down to 109..114).
$ws.Rows("108").Insert()

# Copy formatting (borders/fonts/number format) from the row above (107)
# into the freshly inserted row so it matches the rest of the table.
$ws.Range("A107:D107").Copy($ws.Range("A108:D108"))

# Fill in the new synonym row's content.
$ws.Range("A108").Value2 = "SAVE"
$ws.Range("B108").Value2 = "SAVEIC"
$ws.Range("C108").Value2 = "Save current situation as IC"

# Column C normally doesn't wrap, but on this new row it does (matches
# column B's formatting) - set explicitly to mirror the authored diff.
$ws.Range("C108").WrapText = $true

# Update selection to the cell the author left selected.
$ws.Range("C108").Select()
